$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: "Added in IANA blackhole servers"
# Appends three new rows (78-80) describing the IANA blackhole / prisoner
# servers, each linking Notes -> Wikipedia's Blackhole server article.

$wikiUrl = "https://en.wikipedia.org/wiki/Blackhole_server"

$rows = @(
    @{ Row = 78; Org = "IANA"; IP = "192.175.48.6";  Host = "blackhole-1.iana.org" },
    @{ Row = 79; Org = "IANA"; IP = "192.175.48.42"; Host = "blackhole-2.iana.org" },
    @{ Row = 80; Org = "IANA"; IP = "192.175.48.1";  Host = "prisoner.iana.org" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Range("A$rowNum").Value = $r.Org
    $ws.Range("B$rowNum").Value = $r.IP
    $ws.Range("C$rowNum").Value = $r.Host
    $ws.Range("D$rowNum").Value = $wikiUrl

    $dCell = $ws.Range("D$rowNum")
    $ws.Hyperlinks.Add($dCell, $wikiUrl) | Out-Null
    $dCell.Style = "Hyperlink"
}
